$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    9  = 7
    16 = -6
    20 = -12
    22 = -4
    23 = -3
    24 = -6
    26 = -10
    29 = -2
    30 = -2
    32 = -3
    33 = 4
    36 = -4
    37 = 3
    40 = -3
    42 = 4
    45 = -4
    51 = -2
    53 = 1
    55 = -3
    56 = 7
    58 = -2
    60 = -5
    61 = -1
    66 = 3
    72 = -1
    73 = 5
    76 = 0
    79 = -2
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
